$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last two existing rows were both numbered "25" by mistake; renumber
# them to continue the sequence (26, 27) now that more rows follow.
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27

# --- New activity row: Ansys Fluent simulation (2019.04.09, #meeting) ---
$ws.Range("A29").Value = 28
# Column B stores the date as plain text (like the rest of the sheet), not
# as an Excel date serial. Prefix with an apostrophe so it is entered as
# text, then reset the cell style back to Normal so no visible "quote
# prefix" formatting/indicator is left behind on the cell.
$ws.Range("B29").Value = "'2019.04.09"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "#meeting"
$ws.Range("D29").Value = 2.5
$ws.Range("E29").Value = "Ansys Fluent simulation"

# --- New activity row: Added numerics to equation (2019.04.14, #latex) ---
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "'2019.04.14"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "#latex"
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = "Added numerics to equation for better following"

# Trailing row with just the running number, no activity yet.
$ws.Range("A31").Value = 30

# Column E needs to be widened so the new, longer comment text still fits.
$ws.Columns.Item(5).ColumnWidth = 45.140625

# Leave the selection where it would be after typing the last entry -
# the first empty cell of the next new row.
$ws.Range("A32").Select()
